# Aggiornamento dati: aggiunge le righe dal 2/09/2021 al 9/09/2021 (compreso)
# in fondo alla tabella, replicando formato/valori delle colonne A:D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nuove righe: data (seriale), nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @(44441, 0, 7, 583.8198498748957),
    @(44442, 0, 3, 250.208507089241),
    @(44443, 0, 1, 83.40283569641367),
    @(44444, 0, 0, 0),
    @(44445, 0, 0, 0),
    @(44446, 0, 0, 0),
    @(44447, 0, 0, 0),
    @(44448, 0, 0, 0)
)

$lastRow = 366
$srcRange = $ws.Range("A$lastRow`:D$lastRow")

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $lastRow + 1 + $i
    $destRange = $ws.Range("A$r`:D$r")

    # Copia la formattazione della riga precedente (stile data, bordi, ecc.)
    $srcRange.Copy()
    $destRange.PasteSpecial(-4122)  # xlPasteFormats

    $values = $newRows[$i]
    $ws.Cells.Item($r, 1).Value2 = $values[0]
    $ws.Cells.Item($r, 2).Value2 = $values[1]
    $ws.Cells.Item($r, 3).Value2 = $values[2]
    $ws.Cells.Item($r, 4).Value2 = $values[3]
}

$excel.CutCopyMode = 0
